# Rebuild the 'Article Searches' worksheet: headers + one article row, with a
# hyperlink in the URL column and wrapped/centered formatting, matching the
# target OOXML produced by Excel after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# String literals (PowerShell single-quoted here-strings: no escaping needed,
# content is taken verbatim -- important for the embedded quotes in the search
# term and the newlines in the title).
# ---------------------------------------------------------------------------
$hSearchTerm = @'
Search Term
'@

$hTitle = @'
Title
'@

$hDatabase = @'
Database
'@

$hPubYear = @'
Pub. Year
'@

$hJournal = @'
Journal
'@

$hFirstAuthor = @'
First Author 
'@

$hUrl = @'
URL
'@

$hAbstract = @'
Abstract
'@


$vTitle = @'
Blockchain-Enabled Federated learning
for Enhanced Collaborative Intrusion Detection
in Vehicular Edge Computing
'@

$vAuthor = @'
Zakaria Abou El Houda
'@

$vDatabase = @'
IEEE Xplore
'@

$vSearchTerm = @'
("All Metadata":ai or artificial intelligence) AND ("All Metadata":federated learning)
'@

$vJournal = @'
IEEE Transactions on Intelligent Transportation Systems
'@

$vLink = @'
Link
'@

$vAbstract = @'
Intelligent Transportation Systems (ITSs) are transforming the global monitoring of road safety. These systems, including vehicular networks and transportation infrastructure, are vulnerable to several security issues, which could disrupt services and potentially cause harm to the users. It is crucial to establish robust security measures to protect against evolving attacks and ensure the safe and reliable operation of ITS. Artificial Intelligence (AI)-based Intrusion Detection Systems (IDS) are mainly used to enhance the security of ITS. The adoption of AI-based techniques to secure ITS against new emerging threats has been limited due to a lack of realistic and recent data on these types of attacks ( i.e., zero-day attacks). In this context, we introduce a novel Edge-based Framework that uses Federated Learning (FL) and blockchain to secure ITS against new emerging threats. In particular, our proposed framework consists of a novel distributed Edge-based architecture that allows multiple Edge nodes to securely collaborate while preserving their privacy; and (2) a decentralized and secure reputation system based on blockchain technology to maintain the reliability and trustworthiness of the FL process within the ITS; This system manages reputation data for individual nodes (such as vehicles), guaranteeing the integrity of the FL training process. Experiment results using the UNSW-NB15 dataset show that our proposed framework achieves high accuracy and F1 score (99%) in detecting new threats while ensuring the privacy and reliability of the whole ITS. These results demonstrate the effectiveness of our proposed framework in securing ITS.
'@


$hyperlinkAddress = "https://ieeexplore.ieee.org/document/10542663"

# ---------------------------------------------------------------------------
# Row 1: column headers
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = $hTitle
$ws.Cells.Item(1,2).Value = $hFirstAuthor
$ws.Cells.Item(1,3).Value = $hDatabase
$ws.Cells.Item(1,4).Value = $hSearchTerm
$ws.Cells.Item(1,5).Value = $hPubYear
$ws.Cells.Item(1,6).Value = $hJournal
$ws.Cells.Item(1,7).Value = $hUrl
$ws.Cells.Item(1,8).Value = $hAbstract

# ---------------------------------------------------------------------------
# Row 2: the article
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = $vTitle
$ws.Cells.Item(2,2).Value = $vAuthor
$ws.Cells.Item(2,3).Value = $vDatabase
$ws.Cells.Item(2,4).Value = $vSearchTerm
$ws.Cells.Item(2,5).Value = 2024
$ws.Cells.Item(2,6).Value = $vJournal
$ws.Cells.Item(2,7).Value = $vLink
$ws.Cells.Item(2,8).Value = $vAbstract

# Hyperlink on the URL cell (G2), displaying the text "Link"
[void]$ws.Hyperlinks.Add($ws.Cells.Item(2,7), $hyperlinkAddress, "", "", $vLink)

# ---------------------------------------------------------------------------
# Formatting: the whole used range gets a smaller, centered, word-wrapped
# font; row 2 is tall enough to show the wrapped abstract in full.
# ---------------------------------------------------------------------------
$used = $ws.Range("A1:H2")
$usedFont = $used.Font
$usedFont.Name = "Calibri"
$usedFont.Size = 10
$used.HorizontalAlignment = -4108
$used.VerticalAlignment = -4108
$used.WrapText = $true

$ws.Rows.Item(2).RowHeight = 178.5

# Column widths (character units); the stored OOXML width is
# chars + 5/6, so back the COM input off by 5/6 to land on target widths.
$ws.Columns.Item(1).ColumnWidth = 44.307291666666664
$ws.Columns.Item(2).ColumnWidth = 20.022135416666668
$ws.Columns.Item(3).ColumnWidth = 10.022135416666666
$ws.Columns.Item(4).ColumnWidth = 32.022135416666664
$ws.Columns.Item(5).ColumnWidth = 8.451822916666666
$ws.Columns.Item(6).ColumnWidth = 28.592447916666668
$ws.Columns.Item(8).ColumnWidth = 99.30729166666667

# Active selection, matching the saved view state
[void]$ws.Range("F8").Select()

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

